$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accuracy values for columns C (pre-hardsplit) and D (without pre-hardsplit),
# rows 4-13 (Binkley dataset, runs 1-10). Plain font, no border - matches the
# existing column default formatting.
$ws.Range("C4").Value  = 0.75596815
$ws.Range("D4").Value  = 0.75596815
$ws.Range("C5").Value  = 0.76373625
$ws.Range("D5").Value  = 0.75274724
$ws.Range("C6").Value  = 0.76216215
$ws.Range("D6").Value  = 0.772973
$ws.Range("C7").Value  = 0.7394958
$ws.Range("D7").Value  = 0.7535014
$ws.Range("C8").Value  = 0.7723577
$ws.Range("D8").Value  = 0.7723577
$ws.Range("C9").Value  = 0.7574124
$ws.Range("D9").Value  = 0.7412399
$ws.Range("C10").Value = 0.76151764
$ws.Range("D10").Value = 0.7506775
$ws.Range("C11").Value = 0.75136614
$ws.Range("D11").Value = 0.7021858
$ws.Range("C12").Value = 0.78571427
$ws.Range("D12").Value = 0.771978
$ws.Range("C13").Value = 0.75842696
$ws.Range("D13").Value = 0.7724719

# Row 14 "Average" row: add the AVERAGE formulas for the two new columns,
# matching the bordered / centered look already used by the other Average
# cells on that row (E14, F14).
$ws.Range("C14").Formula = "=AVERAGE(C4:C13)"
$ws.Range("D14").Formula = "=AVERAGE(D4:D13)"
$ws.Range("C14:D14").Borders.LineStyle = 1
$ws.Range("C14:D14").HorizontalAlignment = -4108
$ws.Range("C14:D14").VerticalAlignment = -4108

# Match the recorded selection left behind by the editing session.
$ws.Range("L17").Select()
